$wb = $excel.ActiveWorkbook

# Add new sheet after Sheet2
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "UniqueSheetInFile2"
$newSheet.Range("A1").Value = "Only in File2"
